{"js": "// Applies the relaatiokaavio.docx edits:\n//  1) \"Tehtava(\" row: \"tekija\" -> \"tekija_id\"\n//  2) \"Sessio_tehtava(\" row: add a new \"vastaus_oikein\" field before the\n//     closing parenthesis\n//  3) The \"_GoBack\" bookmark (Word's \"last edit position\" marker) moves\n//     from the very start of the document to right after the newly\n//     typed \"tekija_id\" text, matching where the author's cursor ended\n//     up after the edit.\n\nconst body = context.document.body;\n\n// --- 1) Kayttaja/Tehtava rows: move the _GoBack bookmark & fix \"tekija\" ---\n\n// Remove the bookmark from its original spot (start of the \"Kayttaja(\" row).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// \"tekija\" only occurs once in the document (inside the \"Tehtava(\" row).\nconst tekijaResults = body.search(\"tekija\", { matchCase: true, matchWholeWord: false });\ntekijaResults.load(\"items\");\nawait context.sync();\n\ntekijaResults.items[0].insertText(\"tekija_id\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-insert the bookmark immediately after the text we just typed, i.e.\n// right before \", kuvaus, kyselytyyppi, pvm, esimvastaus)\".\nconst tekijaIdResults = body.search(\"tekija_id\", { matchCase: true, matchWholeWord: false });\ntekijaIdResults.load(\"items\");\nawait context.sync();\n\nconst afterTekijaId = tekijaIdResults.items[0].getRange(Word.RangeLocation.after);\nafterTekijaId.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2) Sessio_tehtava(...) row: add the vastaus_oikein field ---\n\nconst sessioTehtavaResults = body.search(\"aloitus, lopetus)\", { matchCase: true, matchWholeWord: false });\nsessioTehtavaResults.load(\"items\");\nawait context.sync();\n\nsessioTehtavaResults.items[0].insertText(\"aloitus, lopetus, vastaus_oikein)\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Applies the relaatiokaavio.docx edits:\n#  1) \"Tehtava(\" row: \"tekija\" -> \"tekija_id\"\n#  2) \"Sessio_tehtava(\" row: add a new \"vastaus_oikein\" field before the\n#     closing parenthesis\n#  3) The \"_GoBack\" bookmark (Word's \"last edit position\" marker) moves\n#     from the very start of the document to right after the newly\n#     typed \"tekija_id\" text, matching where the author's cursor ended\n#     up after the edit.\n\n$d = $word.ActiveDocument\n\n# --- 1) Kayttaja/Tehtava rows: move the _GoBack bookmark & fix \"tekija\" ---\n\n# Remove the bookmark from its original spot (start of the \"Kayttaja(\" row).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# \"tekija\" only occurs once in the document (inside the \"Tehtava(\" row).\n$range = $d.Content\n$range.Find.Execute(\"tekija\") | Out-Null\n$range.Text = \"tekija_id\"\n\n# Re-insert the bookmark immediately after the text we just typed, i.e.\n# right before \", kuvaus, kyselytyyppi, pvm, esimvastaus)\".\n$bmRange = $d.Range($range.End, $range.End)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# --- 2) Sessio_tehtava(...) row: add the vastaus_oikein field ---\n\n$range2 = $d.Content\n$range2.Find.Execute(\"aloitus, lopetus)\") | Out-Null\n$range2.Text = \"aloitus, lopetus, vastaus_oikein)\"\n"}
